$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 63: hours worked on that day increased from 3 to 5
$ws.Range("C63").Value = 5

# New row 64: Federico Speroni, 2017-06-08, 3 hours, Sprint 3 - Testing / Testing WebAPI
# Copy B63 first so the new date cell (B64) inherits the existing date style (s="1")
# instead of the engine minting a brand-new number-format style.
$ws.Range("B63").Copy($ws.Range("B64"))

$ws.Range("A64").Value = "Federico Speroni"
$ws.Range("B64").Value = 42894
$ws.Range("C64").Value = 3
$ws.Range("D64").Value = "Sprint 3 - Testing"
$ws.Range("E64").Value = "Testing WebAPI"

# Cursor ends up parked on D65 (empty) after entering the row above; a stray
# font touch on that cell leaves a formatting-only style behind (no value).
$ws.Range("D65").Font.ThemeColor = 1
$ws.Range("D65").Select()

# Page setup gets touched (explicit portrait orientation written out)
$ws.PageSetup.Orientation = 1
